# Auto-generated Excel COM-interop script
# Applies updated TPM-derived NATMI values to columns G:T for rows 2-26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1904096666666667
$ws.Range("H2").Value = 0.571229
$ws.Range("I2").Value = 0.09975479936454949
$ws.Range("J2").Value = 0.1077412252037539
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07487166666666667
$ws.Range("N2").Value = 0.224615
$ws.Range("O2").Value = 0.01287435003490057
$ws.Range("P2").Value = 0.01655871537719798
$ws.Range("Q2").Value = 0.01425628909277778
$ws.Range("R2").Value = 0.128306601835
$ws.Range("S2").Value = 0.001284278204680487
$ws.Range("T2").Value = 0.00178405628253955

# Row 3
$ws.Range("G3").Value = 0.1904096666666667
$ws.Range("H3").Value = 0.571229
$ws.Range("I3").Value = 0.09975479936454949
$ws.Range("J3").Value = 0.1077412252037539
$ws.Range("O3").Value = 0.02600892111095355
$ws.Range("P3").Value = 0.03345212152666174
$ws.Range("Q3").Value = 0.02880073148111111
$ws.Range("R3").Value = 0.25920658333
$ws.Range("S3").Value = 0.002594514707111567
$ws.Range("T3").Value = 0.003604172558947406

# Row 4
$ws.Range("G4").Value = 0.1904096666666667
$ws.Range("H4").Value = 0.571229
$ws.Range("I4").Value = 0.09975479936454949
$ws.Range("J4").Value = 0.1077412252037539
$ws.Range("M4").Value = 0.8000470000000001
$ws.Range("N4").Value = 2.400141
$ws.Range("O4").Value = 0.1375698656239178
$ws.Range("P4").Value = 0.1769394371887155
$ws.Range("Q4").Value = 0.1523366825876667
$ws.Range("R4").Value = 1.371030143289
$ws.Range("S4").Value = 0.01372325434392195
$ws.Range("T4").Value = 0.01906367174957486

# Row 5
$ws.Range("G5").Value = 0.1904096666666667
$ws.Range("H5").Value = 0.571229
$ws.Range("I5").Value = 0.09975479936454949
$ws.Range("J5").Value = 0.1077412252037539
$ws.Range("M5").Value = 3.881946
$ws.Range("N5").Value = 7.763892
$ws.Range("O5").Value = 0.6675092708044715
$ws.Range("P5").Value = 0.5723574910282232
$ws.Range("Q5").Value = 0.7391600438780001
$ws.Range("R5").Value = 4.434960263268
$ws.Range("S5").Value = 0.06658725338307679
$ws.Range("T5").Value = 0.06166649733792735

# Row 6
$ws.Range("G6").Value = 0.1904096666666667
$ws.Range("H6").Value = 0.571229
$ws.Range("I6").Value = 0.09975479936454949
$ws.Range("J6").Value = 0.1077412252037539
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9074473333333333
$ws.Range("N6").Value = 2.722342
$ws.Range("O6").Value = 0.1560375924257564
$ws.Range("P6").Value = 0.2006922348792017
$ws.Range("Q6").Value = 0.1727867442575556
$ws.Range("R6").Value = 1.555080698318
$ws.Range("S6").Value = 0.01556549872575868
$ws.Range("T6").Value = 0.02162282727476474

# Row 7
$ws.Range("I7").Value = 0.5298350724050007
$ws.Range("J7").Value = 0.5722539689365677
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07487166666666667
$ws.Range("N7").Value = 0.224615
$ws.Range("O7").Value = 0.01287435003490057
$ws.Range("P7").Value = 0.01655871537719798
$ws.Range("Q7").Value = 0.07572048675166666
$ws.Range("R7").Value = 0.6814843807650001
$ws.Range("S7").Value = 0.006821282182908867
$ws.Range("T7").Value = 0.009475790595092516

# Row 8
$ws.Range("I8").Value = 0.5298350724050007
$ws.Range("J8").Value = 0.5722539689365677
$ws.Range("O8").Value = 0.02600892111095355
$ws.Range("P8").Value = 0.03345212152666174
$ws.Range("S8").Value = 0.01378043859999803
$ws.Range("T8").Value = 0.01914310931298057

# Row 9
$ws.Range("I9").Value = 0.5298350724050007
$ws.Range("J9").Value = 0.5722539689365677
$ws.Range("M9").Value = 0.8000470000000001
$ws.Range("N9").Value = 2.400141
$ws.Range("O9").Value = 0.1375698656239178
$ws.Range("P9").Value = 0.1769394371887155
$ws.Range("Q9").Value = 0.809117132839
$ws.Range("R9").Value = 7.282054195551
$ws.Range("S9").Value = 0.07288933971359469
$ws.Range("T9").Value = 0.1012542951926449

# Row 10
$ws.Range("I10").Value = 0.5298350724050007
$ws.Range("J10").Value = 0.5722539689365677
$ws.Range("M10").Value = 3.881946
$ws.Range("N10").Value = 7.763892
$ws.Range("O10").Value = 0.6675092708044715
$ws.Range("P10").Value = 0.5723574910282232
$ws.Range("Q10").Value = 3.925955621802
$ws.Range("R10").Value = 23.555733730812
$ws.Range("S10").Value = 0.3536698228276964
$ws.Range("T10").Value = 0.3275338458914767

# Row 11
$ws.Range("I11").Value = 0.5298350724050007
$ws.Range("J11").Value = 0.5722539689365677
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.9074473333333333
$ws.Range("N11").Value = 2.722342
$ws.Range("O11").Value = 0.1560375924257564
$ws.Range("P11").Value = 0.2006922348792017
$ws.Range("Q11").Value = 0.9177350637513332
$ws.Range("R11").Value = 8.259615573762
$ws.Range("S11").Value = 0.08267418908080265
$ws.Range("T11").Value = 0.114846927944373

# Row 12
$ws.Range("G12").Value = 0.08741566666666667
$ws.Range("H12").Value = 0.262247
$ws.Range("I12").Value = 0.04579668901430952
$ws.Range("J12").Value = 0.04946319792238989
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.07487166666666667
$ws.Range("N12").Value = 0.224615
$ws.Range("O12").Value = 0.01287435003490057
$ws.Range("P12").Value = 0.01655871537719798
$ws.Range("Q12").Value = 0.006544956656111112
$ws.Range("R12").Value = 0.058904609905
$ws.Range("S12").Value = 0.0005896026048097063
$ws.Range("T12").Value = 0.0008190470160428644

# Row 13
$ws.Range("G13").Value = 0.08741566666666667
$ws.Range("H13").Value = 0.262247
$ws.Range("I13").Value = 0.04579668901430952
$ws.Range("J13").Value = 0.04946319792238989
$ws.Range("O13").Value = 0.02600892111095355
$ws.Range("P13").Value = 0.03345212152666174
$ws.Range("Q13").Value = 0.01322220235444445
$ws.Range("R13").Value = 0.11899982119
$ws.Range("S13").Value = 0.001191122471716049
$ws.Range("T13").Value = 0.001654648907997109

# Row 14
$ws.Range("G14").Value = 0.08741566666666667
$ws.Range("H14").Value = 0.262247
$ws.Range("I14").Value = 0.04579668901430952
$ws.Range("J14").Value = 0.04946319792238989
$ws.Range("M14").Value = 0.8000470000000001
$ws.Range("N14").Value = 2.400141
$ws.Range("O14").Value = 0.1375698656239178
$ws.Range("P14").Value = 0.1769394371887155
$ws.Range("Q14").Value = 0.06993664186966668
$ws.Range("R14").Value = 0.6294297768270001
$ws.Range("S14").Value = 0.006300244353718912
$ws.Range("T14").Value = 0.008751990401941707

# Row 15
$ws.Range("G15").Value = 0.08741566666666667
$ws.Range("H15").Value = 0.262247
$ws.Range("I15").Value = 0.04579668901430952
$ws.Range("J15").Value = 0.04946319792238989
$ws.Range("M15").Value = 3.881946
$ws.Range("N15").Value = 7.763892
$ws.Range("O15").Value = 0.6675092708044715
$ws.Range("P15").Value = 0.5723574910282232
$ws.Range("Q15").Value = 0.339342897554
$ws.Range("R15").Value = 2.036057385324
$ws.Range("S15").Value = 0.0305697144892009
$ws.Range("T15").Value = 0.0283106318610915

# Row 16
$ws.Range("G16").Value = 0.08741566666666667
$ws.Range("H16").Value = 0.262247
$ws.Range("I16").Value = 0.04579668901430952
$ws.Range("J16").Value = 0.04946319792238989
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.9074473333333333
$ws.Range("N16").Value = 2.722342
$ws.Range("O16").Value = 0.1560375924257564
$ws.Range("P16").Value = 0.2006922348792017
$ws.Range("Q16").Value = 0.07932511360822223
$ws.Range("R16").Value = 0.7139260224739999
$ws.Range("S16").Value = 0.007146005094863946
$ws.Range("T16").Value = 0.009926879735316712

# Row 17
$ws.Range("G17").Value = 0.42447
$ws.Range("H17").Value = 0.84894
$ws.Range("I17").Value = 0.2223779938672774
$ws.Range("J17").Value = 0.1601211348241683
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.07487166666666667
$ws.Range("N17").Value = 0.224615
$ws.Range("O17").Value = 0.01287435003490057
$ws.Range("P17").Value = 0.01655871537719798
$ws.Range("Q17").Value = 0.03178077635
$ws.Range("R17").Value = 0.1906846581
$ws.Range("S17").Value = 0.002862972133106301
$ws.Range("T17").Value = 0.002651400297427347

# Row 18
$ws.Range("G18").Value = 0.42447
$ws.Range("H18").Value = 0.84894
$ws.Range("I18").Value = 0.2223779938672774
$ws.Range("J18").Value = 0.1601211348241683
$ws.Range("O18").Value = 0.02600892111095355
$ws.Range("P18").Value = 0.03345212152666174
$ws.Range("Q18").Value = 0.06420391730000001
$ws.Range("R18").Value = 0.3852235038
$ws.Range("S18").Value = 0.005783811699306129
$ws.Range("T18").Value = 0.005356391661125068

# Row 19
$ws.Range("G19").Value = 0.42447
$ws.Range("H19").Value = 0.84894
$ws.Range("I19").Value = 0.2223779938672774
$ws.Range("J19").Value = 0.1601211348241683
$ws.Range("M19").Value = 0.8000470000000001
$ws.Range("N19").Value = 2.400141
$ws.Range("O19").Value = 0.1375698656239178
$ws.Range("P19").Value = 0.1769394371887155
$ws.Range("Q19").Value = 0.33959595009
$ws.Range("R19").Value = 2.03757570054
$ws.Range("S19").Value = 0.03059251073403775
$ws.Range("T19").Value = 0.02833174347780678

# Row 20
$ws.Range("G20").Value = 0.42447
$ws.Range("H20").Value = 0.84894
$ws.Range("I20").Value = 0.2223779938672774
$ws.Range("J20").Value = 0.1601211348241683
$ws.Range("M20").Value = 3.881946
$ws.Range("N20").Value = 7.763892
$ws.Range("O20").Value = 0.6675092708044715
$ws.Range("P20").Value = 0.5723574910282232
$ws.Range("Q20").Value = 1.64776961862
$ws.Range("R20").Value = 6.591078474480001
$ws.Range("S20").Value = 0.1484393725293076
$ws.Range("T20").Value = 0.09164653098855285

# Row 21
$ws.Range("G21").Value = 0.42447
$ws.Range("H21").Value = 0.84894
$ws.Range("I21").Value = 0.2223779938672774
$ws.Range("J21").Value = 0.1601211348241683
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.9074473333333333
$ws.Range("N21").Value = 2.722342
$ws.Range("O21").Value = 0.1560375924257564
$ws.Range("P21").Value = 0.2006922348792017
$ws.Range("Q21").Value = 0.38518416958
$ws.Range("R21").Value = 2.31110501748
$ws.Range("S21").Value = 0.03469932677151959
$ws.Range("T21").Value = 0.03213506839925631

# Row 22
$ws.Range("G22").Value = 0.1951446666666667
$ws.Range("H22").Value = 0.585434
$ws.Range("I22").Value = 0.102235445348863
$ws.Range("J22").Value = 0.1104204731131201
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.07487166666666667
$ws.Range("N22").Value = 0.224615
$ws.Range("O22").Value = 0.01287435003490057
$ws.Range("P22").Value = 0.01655871537719798
$ws.Range("Q22").Value = 0.01461080643444444
$ws.Range("R22").Value = 0.13149725791
$ws.Range("S22").Value = 0.00131621490939521
$ws.Range("T22").Value = 0.001828421186095697

# Row 23
$ws.Range("G23").Value = 0.1951446666666667
$ws.Range("H23").Value = 0.585434
$ws.Range("I23").Value = 0.102235445348863
$ws.Range("J23").Value = 0.1104204731131201
$ws.Range("O23").Value = 0.02600892111095355
$ws.Range("P23").Value = 0.03345212152666174
$ws.Range("Q23").Value = 0.02951693179777778
$ws.Range("R23").Value = 0.26565238618
$ws.Range("S23").Value = 0.002659033632821781
$ws.Range("T23").Value = 0.003693799085611578

# Row 24
$ws.Range("G24").Value = 0.1951446666666667
$ws.Range("H24").Value = 0.585434
$ws.Range("I24").Value = 0.102235445348863
$ws.Range("J24").Value = 0.1104204731131201
$ws.Range("M24").Value = 0.8000470000000001
$ws.Range("N24").Value = 2.400141
$ws.Range("O24").Value = 0.1375698656239178
$ws.Range("P24").Value = 0.1769394371887155
$ws.Range("Q24").Value = 0.1561249051326667
$ws.Range("R24").Value = 1.405124146194
$ws.Range("S24").Value = 0.01406451647864447
$ws.Range("T24").Value = 0.01953773636674715

# Row 25
$ws.Range("G25").Value = 0.1951446666666667
$ws.Range("H25").Value = 0.585434
$ws.Range("I25").Value = 0.102235445348863
$ws.Range("J25").Value = 0.1104204731131201
$ws.Range("M25").Value = 3.881946
$ws.Range("N25").Value = 7.763892
$ws.Range("O25").Value = 0.6675092708044715
$ws.Range("P25").Value = 0.5723574910282232
$ws.Range("Q25").Value = 0.757541058188
$ws.Range("R25").Value = 4.545246349128
$ws.Range("S25").Value = 0.06824310757518995
$ws.Range("T25").Value = 0.06319998494917478

# Row 26
$ws.Range("G26").Value = 0.1951446666666667
$ws.Range("H26").Value = 0.585434
$ws.Range("I26").Value = 0.102235445348863
$ws.Range("J26").Value = 0.1104204731131201
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.9074473333333333
$ws.Range("N26").Value = 2.722342
$ws.Range("O26").Value = 0.1560375924257564
$ws.Range("P26").Value = 0.2006922348792017
$ws.Range("Q26").Value = 0.1770835073808889
$ws.Range("R26").Value = 1.593751566428
$ws.Range("S26").Value = 0.01595257275281158
$ws.Range("T26").Value = 0.02162282727476474
